$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 650-673: week ending 2021-04-25 (Reko93-95, order receipts, bank fee, IKEA, ST1 fuel).
# Column A keeps the same date/time display format used throughout the "Datum" column.
$ws.Range("A650:A673").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 650
$ws.Cells.Item(650, 1).Value = 44305
$ws.Cells.Item(650, 2).Value = "Reko93"
$ws.Cells.Item(650, 3).Value = 3011
$ws.Cells.Item(650, 4).Value = "Reko Swish +46708261018"
$ws.Cells.Item(650, 6).Value = 610.71

# Row 651
$ws.Cells.Item(651, 1).Value = 44305
$ws.Cells.Item(651, 2).Value = "Reko93"
$ws.Cells.Item(651, 3).Value = 2611
$ws.Cells.Item(651, 4).Value = "Reko Swish +46708261018"
$ws.Cells.Item(651, 6).Value = 73.29

# Row 652
$ws.Cells.Item(652, 1).Value = 44305
$ws.Cells.Item(652, 2).Value = "Reko93"
$ws.Cells.Item(652, 3).Value = 1930
$ws.Cells.Item(652, 4).Value = "Reko Swish +46708261018"
$ws.Cells.Item(652, 5).Value = 684

# Row 653
$ws.Cells.Item(653, 1).Value = 44306
$ws.Cells.Item(653, 3).Value = 6570
$ws.Cells.Item(653, 4).Value = "Pris banktjänster enligt faktura"
$ws.Cells.Item(653, 5).Value = 123

# Row 654
$ws.Cells.Item(654, 1).Value = 44306
$ws.Cells.Item(654, 4).Value = "Pris banktjänster enligt faktura"
$ws.Cells.Item(654, 5).Value = 0

# Row 655
$ws.Cells.Item(655, 1).Value = 44306
$ws.Cells.Item(655, 3).Value = 1930
$ws.Cells.Item(655, 4).Value = "Pris banktjänster enligt faktura"
$ws.Cells.Item(655, 6).Value = 123

# Row 656
$ws.Cells.Item(656, 1).Value = 44306
$ws.Cells.Item(656, 2).Formula = "'8201150"
$ws.Cells.Item(656, 3).Value = 3011
$ws.Cells.Item(656, 4).Value = "Order 8201150 Swish +46764282407"
$ws.Cells.Item(656, 6).Value = 935.71

# Row 657
$ws.Cells.Item(657, 1).Value = 44306
$ws.Cells.Item(657, 2).Formula = "'8201150"
$ws.Cells.Item(657, 3).Value = 2611
$ws.Cells.Item(657, 4).Value = "Order 8201150 Swish +46764282407"
$ws.Cells.Item(657, 6).Value = 112.29

# Row 658
$ws.Cells.Item(658, 1).Value = 44306
$ws.Cells.Item(658, 2).Formula = "'8201150"
$ws.Cells.Item(658, 3).Value = 1930
$ws.Cells.Item(658, 4).Value = "Order 8201150 Swish +46764282407"
$ws.Cells.Item(658, 5).Value = 1048

# Row 659
$ws.Cells.Item(659, 1).Value = 44306
$ws.Cells.Item(659, 3).Value = 5460
$ws.Cells.Item(659, 4).Value = "IKEA BARKARBY K0135"
$ws.Cells.Item(659, 5).Value = 4988

# Row 660
$ws.Cells.Item(660, 1).Value = 44306
$ws.Cells.Item(660, 3).Value = 2641
$ws.Cells.Item(660, 4).Value = "IKEA BARKARBY K0135"
$ws.Cells.Item(660, 5).Value = 1247

# Row 661
$ws.Cells.Item(661, 1).Value = 44306
$ws.Cells.Item(661, 3).Value = 1930
$ws.Cells.Item(661, 4).Value = "IKEA BARKARBY K0135"
$ws.Cells.Item(661, 6).Value = 6235

# Row 662
$ws.Cells.Item(662, 1).Value = 44308
$ws.Cells.Item(662, 2).Value = "Reko94"
$ws.Cells.Item(662, 3).Value = 3011
$ws.Cells.Item(662, 4).Value = "Reko Swish +46738070018"
$ws.Cells.Item(662, 6).Value = 739.29

# Row 663
$ws.Cells.Item(663, 1).Value = 44308
$ws.Cells.Item(663, 2).Value = "Reko94"
$ws.Cells.Item(663, 3).Value = 2611
$ws.Cells.Item(663, 4).Value = "Reko Swish +46738070018"
$ws.Cells.Item(663, 6).Value = 88.71

# Row 664
$ws.Cells.Item(664, 1).Value = 44308
$ws.Cells.Item(664, 2).Value = "Reko94"
$ws.Cells.Item(664, 3).Value = 1930
$ws.Cells.Item(664, 4).Value = "Reko Swish +46738070018"
$ws.Cells.Item(664, 5).Value = 828

# Row 665
$ws.Cells.Item(665, 1).Value = 44309
$ws.Cells.Item(665, 2).Value = "Reko95"
$ws.Cells.Item(665, 3).Value = 3011
$ws.Cells.Item(665, 4).Value = "Reko Swish +46703533270"
$ws.Cells.Item(665, 6).Value = 115.18

# Row 666
$ws.Cells.Item(666, 1).Value = 44309
$ws.Cells.Item(666, 2).Value = "Reko95"
$ws.Cells.Item(666, 3).Value = 2611
$ws.Cells.Item(666, 4).Value = "Reko Swish +46703533270"
$ws.Cells.Item(666, 6).Value = 13.82

# Row 667
$ws.Cells.Item(667, 1).Value = 44309
$ws.Cells.Item(667, 2).Value = "Reko95"
$ws.Cells.Item(667, 3).Value = 1930
$ws.Cells.Item(667, 4).Value = "Reko Swish +46703533270"
$ws.Cells.Item(667, 5).Value = 129

# Row 668
$ws.Cells.Item(668, 1).Value = 44309
$ws.Cells.Item(668, 3).Value = 5670
$ws.Cells.Item(668, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(668, 5).Value = 795.1

# Row 669
$ws.Cells.Item(669, 1).Value = 44309
$ws.Cells.Item(669, 3).Value = 2641
$ws.Cells.Item(669, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(669, 5).Value = 198.78

# Row 670
$ws.Cells.Item(670, 1).Value = 44309
$ws.Cells.Item(670, 3).Value = 1930
$ws.Cells.Item(670, 4).Value = "ST1 V#LLINGBY K0135"
$ws.Cells.Item(670, 6).Value = 993.88

# Row 671
$ws.Cells.Item(671, 1).Value = 44311
$ws.Cells.Item(671, 2).Formula = "'7251015"
$ws.Cells.Item(671, 3).Value = 3011
$ws.Cells.Item(671, 4).Value = "Order 7251015 Swish +46703564388"
$ws.Cells.Item(671, 6).Value = 1062.5

# Row 672
$ws.Cells.Item(672, 1).Value = 44311
$ws.Cells.Item(672, 2).Formula = "'7251015"
$ws.Cells.Item(672, 3).Value = 2611
$ws.Cells.Item(672, 4).Value = "Order 7251015 Swish +46703564388"
$ws.Cells.Item(672, 6).Value = 127.5

# Row 673
$ws.Cells.Item(673, 1).Value = 44311
$ws.Cells.Item(673, 2).Formula = "'7251015"
$ws.Cells.Item(673, 3).Value = 1930
$ws.Cells.Item(673, 4).Value = "Order 7251015 Swish +46703564388"
$ws.Cells.Item(673, 5).Value = 1190

